$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")
$ws.Activate()

# 1. Update the "Nationally, between 2009 and 2015..." paragraph (row 6, col B)
$ws.Range("B6").Value = "Nationally, between 2009 and 2015, there was no significant change in the proportion of people with disability who report a need for more formal assistance. Progress will need to improve in order to meet the target."

# 2. Insert a new row after row 6 for the new NDIS paragraph.
#    Inserting at row 7 pushes the old rows 7-10 down to 8-11 and copies
#    formatting from the row above (row 6), matching A7 blank / B7 styled like B6.
$ws.Rows.Item(7).Insert()
$ws.Range("B7").Value = "All jurisdictions have committed to rolling out the National Disability Insurance Scheme (NDIS) to full scheme.  The NDIS is transforming Australia" + [char]8217 + "s system of support for people with disability through an insurance approach.  The NDIS provides assurance that people with disability will be able to receive care and support over their lifetime based on their needs, and that they will have choice and control over that support. "

# 3. Old row 9 ("Sourced from: ...") is now row 10 after the insert above.
#    Split it into a "Source" label (A10) and a shortened source string (B10).
$ws.Range("A10").Value = "Source"
$ws.Range("B10").Value = "ABS, Survey of Disability, Ageing and Carers"

# 4. The old Northern Territory caveat row (now row 11) is removed entirely.
$ws.Rows.Item(11).Delete()

# 5. Row heights as re-flowed by the editor for the updated text.
$ws.Rows.Item(5).RowHeight = 25.45
$ws.Rows.Item(6).RowHeight = 25.45
$ws.Rows.Item(7).RowHeight = 49.45
$ws.Rows.Item(8).RowHeight = 20.95
$ws.Rows.Item(9).RowHeight = 20.95
$ws.Rows.Item(10).RowHeight = 12.8
$ws.Rows.Item(11).RowHeight = 13.8

# 6. Leave the selection/cursor where the editor left it when saving.
$wb.Worksheets.Item("Data").Range("A1").Select()
$ws.Activate()
$ws.Range("I26").Select()
